$wb = $excel.ActiveWorkbook
$wsTypography = $wb.Worksheets.Item("Typography")
$wsTranslation = $wb.Worksheets.Item("Translation")

# --- Typography sheet: row 9 (Digital_Dream_25) wildcard/fallback updates ---
$wsTypography.Range("F9").Value2 = "?"
$wsTypography.Range("G9").Value2 = ".-"
$wsTypography.Range("I9").Value2 = "0-9,"

# --- Translation sheet: update existing accelerometer placeholder rows ---
$wsTranslation.Range("F15").Value2 = "X:<value>"
$wsTranslation.Range("F17").Value2 = "Y:<value>"
$wsTranslation.Range("F19").Value2 = "Z:<value>"

# F16/F18/F20 become "-000.0" which Excel would otherwise parse as a number;
# force text storage (matches the shared-string/text cell type in the source)
$wsTranslation.Range("F16").NumberFormat = "@"
$wsTranslation.Range("F16").Value2 = "-000.0"
$wsTranslation.Range("F18").NumberFormat = "@"
$wsTranslation.Range("F18").Value2 = "-000.0"
$wsTranslation.Range("F20").NumberFormat = "@"
$wsTranslation.Range("F20").Value2 = "-000.0"

# --- Translation sheet: new rows 21-24 ---
$wsTranslation.Range("B21").Value2 = "SingleUseId24"
$wsTranslation.Range("C21").Value2 = "Default"
$wsTranslation.Range("D21").Value2 = "Left"
$wsTranslation.Range("E21").Value2 = "LTR"
$wsTranslation.Range("F21").Value2 = "RED"

$wsTranslation.Range("B22").Value2 = "SingleUseId25"
$wsTranslation.Range("C22").Value2 = "Default"
$wsTranslation.Range("D22").Value2 = "Left"
$wsTranslation.Range("E22").Value2 = "LTR"
$wsTranslation.Range("F22").Value2 = "GREEN"

$wsTranslation.Range("B23").Value2 = "SingleUseId26"
$wsTranslation.Range("C23").Value2 = "Large"
$wsTranslation.Range("D23").Value2 = "Left"
$wsTranslation.Range("E23").Value2 = "LTR"
$wsTranslation.Range("F23").Value2 = "LED"

$wsTranslation.Range("B24").Value2 = "SingleUseId27"
$wsTranslation.Range("C24").Value2 = "Digital_Dream_25"
$wsTranslation.Range("D24").Value2 = "Left"
$wsTranslation.Range("E24").Value2 = "LTR"
$wsTranslation.Range("F24").Value2 = "TILT"
